$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New topic "study_books" data for rows 126-129 (mirrors the pattern used by the
# other 4-point topic polygons already in the sheet).
$rows = @(
    @{ r = 126; b = 1; c = 626.172; d = 182.411 },
    @{ r = 127; b = 2; c = 714.949; d = 182.411 },
    @{ r = 128; b = 3; c = 714.949; d = 210.943 },
    @{ r = 129; b = 4; c = 626.172; d = 210.943 }
)

foreach ($row in $rows) {
    $r = $row.r
    $ws.Cells.Item($r, 1).Value = "study_books"
    $ws.Cells.Item($r, 2).Value = $row.b
    $ws.Cells.Item($r, 3).Value = $row.c
    $ws.Cells.Item($r, 4).Value = $row.d
    $ws.Cells.Item($r, 5).Value = 1002.816
    $ws.Cells.Item($r, 6).Value = 225.67
    $ws.Cells.Item($r, 7).Formula = "=C$r/E$r"
    $ws.Cells.Item($r, 8).Formula = "=D$r/F$r"
}

# Restore the view: scroll the frozen pane down to the newly added rows.
$ws.Application.ActiveWindow.ScrollRow = 116
$ws.Range("L123").Select()
